# Update "想去人数" (column F) counts that changed between scrapes.
# Values taken from the commit's regenerated data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    2  = 1550
    4  = 476
    5  = 917
    10 = 5697
    14 = 7939
    15 = 9332
    17 = 927
    18 = 4546
    19 = 694
    20 = 265
    22 = 293
    25 = 127
    26 = 1707
    28 = 974
    32 = 2354
    35 = 1501
    38 = 5
    39 = 809
    40 = 525
    41 = 4177
    42 = 199
    43 = 53
    46 = 16
    48 = 183
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F27").Value = 100

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5368

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    2  = 1550
    4  = 476
    5  = 917
    10 = 5697
    12 = 7939
    13 = 9332
    16 = 927
    17 = 694
    18 = 265
    20 = 293
    24 = 127
    26 = 974
    30 = 2354
    39 = 525
    40 = 4177
    42 = 199
    43 = 53
    46 = 16
    48 = 183
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
